# Rename the unclear month-abbreviation headers (JAN..DEC) on the "Pool"
# sheet to dependents_1..dependents_12 (fixes #13).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pool")

$ws.Range("N1").Value = "dependents_1"
$ws.Range("O1").Value = "dependents_2"
$ws.Range("P1").Value = "dependents_3"
$ws.Range("Q1").Value = "dependents_4"
$ws.Range("R1").Value = "dependents_5"
$ws.Range("S1").Value = "dependents_6"
$ws.Range("T1").Value = "dependents_7"
$ws.Range("U1").Value = "dependents_8"
$ws.Range("V1").Value = "dependents_9"
$ws.Range("W1").Value = "dependents_10"
$ws.Range("X1").Value = "dependents_11"
$ws.Range("Y1").Value = "dependents_12"
